$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 29. This shifts the existing rows 29-55
# down to rows 30-56, preserving all their data/formatting.
$ws.Rows.Item(29).Insert()

# Populate the newly inserted row 29 with the new weekly data point.
$ws.Range("A29").Value = 7
$ws.Range("B29").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C29").Value = 'Ñuble'
$ws.Range("D29").Value = 44566
$ws.Range("E29").Value = 16
$ws.Range("F29").Value = 100112031
$ws.Range("G29").Value = 'Poroto verde'
$ws.Range("H29").Value = 'Sin especificar'
$ws.Range("I29").Value = 'Primera'
$ws.Range("J29").Value = 120
$ws.Range("K29").Value = 38000
$ws.Range("L29").Value = 39000
$ws.Range("M29").Value = 38500
$ws.Range("N29").Value = '$/saco 25 kilos'
$ws.Range("O29").Value = 'Región del Maule'
$ws.Range("P29").Value = 1540
$ws.Range("Q29").Value = 25
$ws.Range("R29").Value = 'Hortaliza'
